# Update the "Förändrad" (changed) date column (C) for all data rows
# from 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C146")
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    $cell = $rng.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
